# Generate Report for handoff
# - Removes the stale "Handoff transform failed" row (old 2ee0c804... entry) from every sheet
# - Refreshes the source-file GUID, the handoff package hash, and the handoff timestamps
#   to the values produced by the latest handoff run

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldGuid = "4caeda85-fa12-4b83-ab12-a5dd74fe18a2"
$newGuid = "93196106-1bf7-426f-b972-12d103c4f071"

$newMdName      = "$newGuid.md"
$newZhXlfName   = "$newGuid.b409b808629d7314fd35d71ec7b987486b60140f.zh-cn.xlf"
$newDeXlfName   = "$newGuid.b409b808629d7314fd35d71ec7b987486b60140f.de-de.xlf"
$newZhDatetime  = "2016-01-15 08:10:35"
$newDeDatetime  = "2016-01-15 08:10:46"

$mdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/18ab818b81a056952a7725b35abbedb1c6964284/e2e/$newMdName"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/18ab818b81a056952a7725b35abbedb1c6964284/.localization-config"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25d6722963afeea392befb06f17054481049e2e9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$newZhXlfName"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e019cd5891d949fe44b8e9adc24ab2dbe72bdfc3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$newDeXlfName"

# --- Overview sheet ---------------------------------------------------
# Row 3 (2ee0c804.../"Handoff transform failed") goes away entirely;
# the ".localization-config" row slides up from row 4 to row 3.
$wsOverview.Rows("3:3").Delete()
$wsOverview.Range("A2").Value = $newMdName

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configUrl, "", "", ".localization-config")

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Rows("3:3").Delete()
$wsZhCn.Range("A2").Value = $newMdName
$wsZhCn.Range("C2").Value = $newZhXlfName
$wsZhCn.Range("D2").Value = $newZhDatetime

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl, "", "", $newMdName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $zhXlfUrl, "", "", $newZhXlfName)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $configUrl, "", "", ".localization-config")

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Rows("3:3").Delete()
$wsDeDe.Range("A2").Value = $newMdName
$wsDeDe.Range("C2").Value = $newDeXlfName
$wsDeDe.Range("D2").Value = $newDeDatetime

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl, "", "", $newMdName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $deXlfUrl, "", "", $newDeXlfName)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $configUrl, "", "", ".localization-config")
